# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamp cells to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-31 07:14:24"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-31 07:14:18"
$zhcn.Range("K2").Value = "2016-08-31 07:14:37"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$dede.Range("H2").Value = "2016-08-31 07:14:24"
$dede.Range("K2").Value = "2016-08-31 07:14:44"
